# "From bjo-sea01 add seamless01" — insert a new USER_LIST row for
# SEAMLESS01 / CLIENT_ADMIN just above the trailing "*END*" marker row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 currently holds the "*END*" sentinel. Insert a blank row above it
# (shifting "*END*" down to row 5) and populate the new row 4.
$ws.Rows("4").Insert()

$ws.Range("A4").Value = "SEAMLESS01"
$ws.Range("B4").Value = "CLIENT_ADMIN"

# Widen the columns to fit the new, longer values (user resized A & B).
$ws.Columns("A").ColumnWidth = 14.8
$ws.Columns("B").ColumnWidth = 17.166666666666668

# Leave the cursor where the author left it.
$ws.Range("B9").Select()
